$wb = $excel.ActiveWorkbook

$wsCustomer = $wb.Worksheets.Item("customer")
$wsAddress  = $wb.Worksheets.Item("address")
$wsPhone    = $wb.Worksheets.Item("phone")

# ---------------------------------------------------------------------------
# Write the new string-valued cells first, in the precise order the source
# data was authored, so the shared-strings table is rebuilt in the same
# sequence as the target workbook.
# ---------------------------------------------------------------------------

# address: street (B) + phone reference (F) for the 5 new rows, row by row
$wsAddress.Range("B4").Value = "street address3"
$wsAddress.Range("F4").Value = "reference:phone@id#2"
$wsAddress.Range("B5").Value = "street address4"
$wsAddress.Range("F5").Value = "reference:phone@id#3"
$wsAddress.Range("B6").Value = "street address5"
$wsAddress.Range("F6").Value = "reference:phone@id#4"
$wsAddress.Range("B7").Value = "street address6"
$wsAddress.Range("F7").Value = "reference:phone@id#5"
$wsAddress.Range("B8").Value = "street address7"
$wsAddress.Range("F8").Value = "reference:phone@id#6"

# customer: address reference (D) for the 5 new rows that only ref one address
$wsCustomer.Range("D3").Value = "listReference:address@id#1"
$wsCustomer.Range("D4").Value = "listReference:address@id#2"
$wsCustomer.Range("D5").Value = "listReference:address@id#3"
$wsCustomer.Range("D6").Value = "listReference:address@id#4"
$wsCustomer.Range("D7").Value = "listReference:address@id#5"

# customer: name (B) for all 6 new rows
$wsCustomer.Range("B3").Value = "customer2"
$wsCustomer.Range("B4").Value = "customer3"
$wsCustomer.Range("B5").Value = "customer4"
$wsCustomer.Range("B6").Value = "customer5"
$wsCustomer.Range("B7").Value = "customer6"
$wsCustomer.Range("B8").Value = "customer7"

# address: city (C) - row 3 changes, rows 4-8 are new
$wsAddress.Range("C3").Value = "New Hamisphere"
$wsAddress.Range("C4").Value = "North Carolina"
$wsAddress.Range("C5").Value = "Chicago"
$wsAddress.Range("C6").Value = "California"
$wsAddress.Range("C7").Value = "Texas"
$wsAddress.Range("C8").Value = "Des Moines"

# customer: last row's address reference is a list of two addresses
$wsCustomer.Range("D8").Value = "listReference:address@id#6,address@id#7"

# ---------------------------------------------------------------------------
# Remaining (numeric / already-known-string) cells - order is immaterial for
# the shared-strings table.
# ---------------------------------------------------------------------------

# phone: new rows 3-7
$phoneRows = @(
    @(2, 123457, 1516),
    @(3, 123458, 1517),
    @(4, 123459, 1518),
    @(5, 123460, 1519),
    @(6, 123461, 1520)
)
$r = 3
foreach ($row in $phoneRows) {
    $wsPhone.Cells.Item($r, 1).Value = $row[0]
    $wsPhone.Cells.Item($r, 2).Value = $row[1]
    $wsPhone.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# address: id (A), zipCode (D), country (E) for new rows 4-8
$addressRows = @(
    @(4, 3, 54322),
    @(5, 4, 54323),
    @(6, 5, 54324),
    @(7, 6, 54325),
    @(8, 7, 54326)
)
foreach ($row in $addressRows) {
    $row1 = $row[0]
    $wsAddress.Cells.Item($row1, 1).Value = $row[1]
    $wsAddress.Cells.Item($row1, 4).Value = $row[2]
    $wsAddress.Cells.Item($row1, 5).Value = "US"
}

# customer: id (A), age (C) for new rows 3-8
$customerRows = @(
    @(3, 2, 31),
    @(4, 3, 32),
    @(5, 4, 33),
    @(6, 5, 34),
    @(7, 6, 35),
    @(8, 7, 36)
)
foreach ($row in $customerRows) {
    $row1 = $row[0]
    $wsCustomer.Cells.Item($row1, 1).Value = $row[1]
    $wsCustomer.Cells.Item($row1, 3).Value = $row[2]
}

# ---------------------------------------------------------------------------
# Column widths (best-effort; engine quantizes to 1/6-character increments)
# ---------------------------------------------------------------------------
$wsCustomer.Columns("D").ColumnWidth = 47.25
$wsAddress.Columns("F").ColumnWidth = 23.59
$wsPhone.Columns("B").ColumnWidth = 10.92
$wsPhone.Columns("C").ColumnWidth = 8.92

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------
$wsPhone.Range("A2:XFD2").Select() | Out-Null
$wsAddress.Range("F8").Select() | Out-Null
$wsCustomer.Range("D8").Select() | Out-Null

$wsCustomer.Activate() | Out-Null
